# Automatische test-sync: 2025-06-24 19:44:50
#
# Adds the new "Offerte voor 500 stuks" mail-log entry (row 6) to the
# "Logs" sheet, adds the matching dashboard aggregate (row 5) to the
# "Dashboard" sheet, extends the conditional formatting ranges on the
# Logs sheet to include the new row, and extends the bar chart's
# category/value series references to include the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 6
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A6").Value = "Offerte voor 500 stuks"
$wsLogs.Range("B6").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C6").Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$wsLogs.Range("D6").Value = "Offerte / Prijsaanvraag"
$wsLogs.Range("E6").Value = "Beste klant,`nDank u voor uw interesse in product X. Om u een nauwkeurige offerte te kunnen sturen, hebben we wat meer details nodig. Kunt u ons laten weten of u specifieke wensen heeft met betrekking tot de productspecificaties of de leveringsvoorwaarden? Zou u ons ook de gewenste leverdatum kunnen doorgeven?`nZodra we deze informatie hebben ontvangen, zullen we zo snel mogelijk een offerte voor 500 stuks product X voor u opstellen.`nMet vriendelijke groet,`n[Naam Bedrijf] E-mailassistent"
$wsLogs.Range("F6").Value = "2025-06-24 19:44:36"
$wsLogs.Range("G6").Value = "Ja"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges to include row 6
#    (D2:D5 -> D2:D6, G2:G5 -> G2:G6)
# ---------------------------------------------------------------------
$fcsD = $wsLogs.Range("D2:D5").FormatConditions()
for ($i = 1; $i -le $fcsD.Count(); $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D6"))
}

$fcsG = $wsLogs.Range("G2:G5").FormatConditions()
for ($i = 1; $i -le $fcsG.Count(); $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G6"))
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append row 5 with the new category aggregate
# ---------------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("A5").Value = "Offerte / Prijsaanvraag"
$wsDash.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard chart: extend the series category/value references
#    ($A$2:$A$4 -> $A$2:$A$5, $B$2:$B$4 -> $B$2:$B$5)
# ---------------------------------------------------------------------
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$5,Dashboard!`$B`$2:`$B`$5,1)"
